$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AR2:AR51").Value = 14
